$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# --- Correct existing histogram delimiter values in row 3 ---
$ws.Range("C3").Value = 460
$ws.Range("D3").Value = 2596

# --- Add new "grids" pile-up correction block (F3:H3) ---
$ws.Range("F3").Value = "grids"
$ws.Range("G3").Value = 460
$ws.Range("H3").Value = 2600

# --- Add new J:K delimiter table (pile-up correction) ---
$ws.Range("J5").Value = 460
$ws.Range("K5").Value = 786

$ws.Range("J6").Value = 810
$ws.Range("K6").Value = 1140

$ws.Range("J7").Value = 1160
$ws.Range("K7").Value = 1486

$ws.Range("J8").Value = 1510
$ws.Range("K8").Value = 1850

$ws.Range("J9").Value = 1882
$ws.Range("K9").Value = 2218

$ws.Range("J10").Value = 2238
$ws.Range("K10").Value = 2575

# --- Update active selection to match the saved view state ---
$ws.Range("D3").Select()
